# Atualização de bases das ligas, do dia: 24-02-2024 às 12:40
#
# The match-result columns (B, F:AC) for a handful of rows were shuffled
# between rows that share the same date - i.e. the row's sequence number
# (col A), division columns (C, D) and date (col E) stay put, while the
# actual match data (id, teams, score, odds, ...) moves to a different row.
#
# Pairs that fully swap their B:AC content:
#   24 <-> 25
#   26 <-> 27
#   28 <-> 29
#   30 <-> 32
#
# A 3-way rotation of B:AC content:
#   95 -> 96 -> 97 -> 95   (new95 = old96, new96 = old97, new97 = old95)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

# Full-row swaps
Swap-Rows 24 25
Swap-Rows 26 27
Swap-Rows 28 29
Swap-Rows 30 32

# 3-way rotation: 95 <- 96 <- 97 <- 95
$r95 = Get-RowValues 95
$r96 = Get-RowValues 96
$r97 = Get-RowValues 97

Set-RowValues 95 $r96
Set-RowValues 96 $r97
Set-RowValues 97 $r95
